# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.200.90"
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").Value = "3.713.71"
$ws.Range("E3").Value = "  +1.06%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'238.28"
$ws.Range("E5").Value = "  -0.60%  "

$ws.Range("D6").Value = "'1.90"
$ws.Range("E6").Value = "  +1.99%  "

$ws.Range("D7").Value = "'662.09"
$ws.Range("E7").Value = "  +0.97%  "

$ws.Range("D8").Value = "'0.426"
$ws.Range("E8").Value = "  +1.22%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("E10").Value = "  -1.65%  "

$ws.Range("D11").Value = "3.711.81"
$ws.Range("E11").Value = "  +1.03%  "

$ws.Range("D12").Value = "'0.0000328"
$ws.Range("E12").Value = "  +22.31%  "

$ws.Range("D13").Value = "'44.53"
$ws.Range("E13").Value = "  -2.33%  "

$ws.Range("E14").Value = "  +1.62%  "

$ws.Range("E15").Value = "  +1.27%  "

$ws.Range("D16").Value = "4.407.57"
$ws.Range("E16").Value = "  +1.16%  "

$ws.Range("D17").Value = "96.997.28"
$ws.Range("E17").Value = "  +0.67%  "

$ws.Range("D18").Value = "'9.15"
$ws.Range("E18").Value = "  +17.77%  "

$ws.Range("D19").Value = "3.709.52"
$ws.Range("E19").Value = "  +0.99%  "

$ws.Range("D20").Value = "'13.07"
$ws.Range("E20").Value = "  +2.07%  "

$ws.Range("D21").Value = "'18.86"
$ws.Range("E21").Value = "  +0.49%  "

$ws.Range("D22").Value = "'0.504"
$ws.Range("E22").Value = "  -3.87%  "

$ws.Range("D23").Value = "'525.85"
$ws.Range("E23").Value = "  -0.72%  "

$ws.Range("D24").Value = "'3.46"
$ws.Range("E24").Value = "  +0.49%  "

$ws.Range("D25").Value = "'0.0000221"
$ws.Range("E25").Value = "  +8.52%  "

$ws.Range("E26").Value = "  -2.88%  "

$ws.Range("D27").Value = "'106.84"
$ws.Range("E27").Value = "  +5.23%  "

$ws.Range("E28").Value = "  +15.62%  "

$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'13.67"
$ws.Range("E29").Value = "  +3.68%  "

$ws.Range("B30").Value = "WrappedeETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D30").Value = "3.913.02"
$ws.Range("E30").Value = "  +1.16%  "

$ws.Range("D31").Value = "'13.08"
$ws.Range("E31").Value = "  +4.97%  "

$ws.Range("D32").Value = "'3.06"
$ws.Range("E32").Value = "  +0.63%  "

$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("E34").Value = "  +4.33%  "

$ws.Range("E35").Value = "  -1.65%  "

$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("D37").Value = "'32.45"
$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("D38").Value = "'650.82"
$ws.Range("E38").Value = "  -3.11%  "

$ws.Range("E39").Value = "  +1.06%  "

$ws.Range("E40").Value = "  +0.43%  "

$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("E42").Value = "  +4.58%  "

$ws.Range("D43").Value = "'6.90"
$ws.Range("E43").Value = "  +6.50%  "

$ws.Range("D44").Value = "'2.05"
$ws.Range("E44").Value = "  +3.49%  "

$ws.Range("D45").Value = "'0.486"
$ws.Range("E45").Value = "  +8.43%  "

$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "'0.977"
$ws.Range("E46").Value = "  +1.96%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'39.77"
$ws.Range("E47").Value = "  +2.37%  "

$ws.Range("D48").Value = "'0.0464"
$ws.Range("E48").Value = "  -0.50%  "

$ws.Range("E49").Value = "  +5.88%  "

$ws.Range("D50").Value = "'8.83"
$ws.Range("E50").Value = "  +3.15%  "

$ws.Range("D51").Value = "'23.62"
$ws.Range("E51").Value = "  -0.16%  "
